$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Split the old combined "ventral_right" / "ventral_left" / "dorsal_right" /
# "dorsal_left" label (column C) into two separate columns: Hemi
# (ventral/dorsal) and ROI_Location (right/left), then run a couple of
# T-Test comparisons across the new groupings.
#
# (Writes are ordered to match the author's original entry order so new
# shared-string values land in the same slots as the source workbook.)

# ventral_right block (rows 2-13): Hemi = ventral, ROI_Location = right
$ws.Range("G2:G13").Value = "ventral"

# New header cells for the added "Hemi" / "ROI_Location" breakdown columns
$ws.Range("G1").Value = "Hemi"
$ws.Range("H1").Value = "ROI_Location"

# ventral_left block (rows 14-25): Hemi = ventral, ROI_Location = left
$ws.Range("G14:G25").Value = "ventral"

# dorsal_right block (rows 26-37): Hemi = dorsal, ROI_Location = right
$ws.Range("G26:G37").Value = "dorsal"

# dorsal_left block (rows 38-49): Hemi = dorsal, ROI_Location = left
$ws.Range("G38:G49").Value = "dorsal"

$ws.Range("H2:H13").Value = "right"
$ws.Range("H14:H25").Value = "left"
$ws.Range("H26:H37").Value = "right"
$ws.Range("H38:H49").Value = "left"

# Additional t-test analysis comparing the slope sensitivity across groups
$ws.Range("I26").Formula = "=T.TEST(B2:B13,B14:B25,2,1)"
$ws.Range("I27").Formula = "=T.TEST(B26:B37,B38:B49,2,1)"
$ws.Range("I28").Formula = "=T.TEST(B2:B25,B26:B49,2,1)"

# Restore the author's last active selection
$ws.Range("L32").Select()
